$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Diebold-Mariano statistics (DM_Stat, column C) and p-values (P_Value, column D)
$updates = @(
    @{ Row = 2;  C = -0.8411854917792635; D = 0.4061215012744057 }
    @{ Row = 3;  C = -0.2778016380342694; D = 0.7828469691540081 }
    @{ Row = 4;  C = -1.518784617224351;  D = 0.1380617964729218 }
    @{ Row = 5;  C = -0.3380130701987242; D = 0.7374316209005798 }
    @{ Row = 6;  C = 0.5045919890224879;  D = 0.6171002573772766 }
    @{ Row = 7;  C = -0.3503126591275078; D = 0.7282651074778745 }
    @{ Row = 8;  C = 0.2863567500055977;  D = 0.7763425859844819 }
    @{ Row = 9;  C = -0.9466563850380291; D = 0.3504937734459981 }
    @{ Row = 10; C = -0.1056209815844399; D = 0.9165037705244508 }
    @{ Row = 11; C = 0.9713308383160291;  D = 0.33824342429289 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 4).Value = $u.D
}
